# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet
# (rows 2-48) to reflect the latest scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    "2" = "2025-10-31T07:01:50.533873+00:00"
    "3" = "2025-10-31T07:01:53.208599+00:00"
    "4" = "2025-10-31T07:01:53.208630+00:00"
    "5" = "2025-10-31T07:01:53.208650+00:00"
    "6" = "2025-10-31T07:01:53.208667+00:00"
    "7" = "2025-10-31T07:01:55.564044+00:00"
    "8" = "2025-10-31T07:01:55.564076+00:00"
    "9" = "2025-10-31T07:01:57.864150+00:00"
    "10" = "2025-10-31T07:02:00.644788+00:00"
    "11" = "2025-10-31T07:02:00.644817+00:00"
    "12" = "2025-10-31T07:02:03.507666+00:00"
    "13" = "2025-10-31T07:02:03.507694+00:00"
    "14" = "2025-10-31T07:02:03.507713+00:00"
    "15" = "2025-10-31T07:02:11.051218+00:00"
    "16" = "2025-10-31T07:02:13.408140+00:00"
    "17" = "2025-10-31T07:02:15.752756+00:00"
    "18" = "2025-10-31T07:02:18.050314+00:00"
    "19" = "2025-10-31T07:02:18.050344+00:00"
    "20" = "2025-10-31T07:02:18.050362+00:00"
    "21" = "2025-10-31T07:02:20.433948+00:00"
    "22" = "2025-10-31T07:02:20.433977+00:00"
    "23" = "2025-10-31T07:02:20.433995+00:00"
    "24" = "2025-10-31T07:02:20.434011+00:00"
    "25" = "2025-10-31T07:02:20.434028+00:00"
    "26" = "2025-10-31T07:02:28.316231+00:00"
    "27" = "2025-10-31T07:02:28.316262+00:00"
    "28" = "2025-10-31T07:02:28.316281+00:00"
    "29" = "2025-10-31T07:02:28.316299+00:00"
    "30" = "2025-10-31T07:02:31.125876+00:00"
    "31" = "2025-10-31T07:02:31.125903+00:00"
    "32" = "2025-10-31T07:02:31.125920+00:00"
    "33" = "2025-10-31T07:02:33.460988+00:00"
    "34" = "2025-10-31T07:02:33.461016+00:00"
    "35" = "2025-10-31T07:02:33.461034+00:00"
    "36" = "2025-10-31T07:02:33.461050+00:00"
    "37" = "2025-10-31T07:02:33.461065+00:00"
    "38" = "2025-10-31T07:02:33.461080+00:00"
    "39" = "2025-10-31T07:02:33.461099+00:00"
    "40" = "2025-10-31T07:02:33.461114+00:00"
    "41" = "2025-10-31T07:02:33.461128+00:00"
    "42" = "2025-10-31T07:02:36.189530+00:00"
    "43" = "2025-10-31T07:02:36.189559+00:00"
    "44" = "2025-10-31T07:02:40.887512+00:00"
    "45" = "2025-10-31T07:02:43.693787+00:00"
    "46" = "2025-10-31T07:02:43.693816+00:00"
    "47" = "2025-10-31T07:02:43.693834+00:00"
    "48" = "2025-10-31T07:02:43.693851+00:00"
}

foreach ($row in $timestamps.Keys) {
    $cell = $ws.Cells.Item([int]$row, 11)  # column K = 11
    $cell.Value = $timestamps[$row]
}
